$wb = $excel.ActiveWorkbook

$newRunId = "4d674e3c0c6d488191ddeab195b0d07a"
$newTimestamp = "2025-10-27T00:03:56.362674"

# --- Sheet "preguntas": rows 2-23 get new run_id/timestamp values ---
$wsPreguntas = $wb.Worksheets.Item("preguntas")

for ($row = 2; $row -le 23; $row++) {
    $wsPreguntas.Range("V$row").Value = $newRunId
    $wsPreguntas.Range("W$row").Value = $newTimestamp
    $wsPreguntas.Range("AC$row").Value = $newRunId
    $wsPreguntas.Range("AD$row").Value = "[`"$newRunId`"]"
}

# --- Sheet "indice_global": row 2 gets new run_id/timestamp/extra_config ---
$wsIndice = $wb.Worksheets.Item("indice_global")

$wsIndice.Range("M2").Value = $newRunId
$wsIndice.Range("Q2").Value = $newTimestamp
$wsIndice.Range("R2").Value = "{'model_name': 'gpt-4o-mini', 'retries': 2, 'backoff_factor': 2.0, 'timeout_seconds': 60.0, 'prompt_batch_size': 1, 'log_level': 'INFO', 'log_file': None, 'ai_provider': 'mock', 'run_id': '$newRunId', 'document_id': None, 'extra_instructions': None, 'splitter_log_level': 'info', 'splitter_normalize_newlines': True}"
